# Quarterly balance-sheet roll-forward:
# drop the oldest reporting quarter (column D) and shift every later quarter
# one column to the left, then append the newest quarter's figures in the
# now-empty last column (M). One shifted-in date label also gets corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift everything from column E..M one column to the left (into D..L),
# which is exactly what dropping the oldest quarter does.
$ws.Columns("D").Delete()

# --- Row 8: quarter labels -------------------------------------------------
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# --- Row 9: publish dates ---------------------------------------------------
# One of the shifted-in labels itself needed correcting (unrelated small fix
# riding along with the roll-forward), and the brand new quarter's date goes
# in M9.
$ws.Range("I9").Value = "1402-02-23 (7)"
$ws.Range("M9").Value = "1402-02-23"

# --- Newest quarter's figures (column M) for every data row ---------------
$ws.Range("M12").Value = 315936
$ws.Range("M13").Value = 366791729
$ws.Range("M14").Value = 108064141
$ws.Range("M15").Value = 104315370
$ws.Range("M16").Value = 1015312
$ws.Range("M17").Value = 0
$ws.Range("M18").Value = 580502488
$ws.Range("M19").Value = 15000
$ws.Range("M20").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("M22").Value = 68104029
$ws.Range("M23").Value = 892305
$ws.Range("M24").Value = "-"
$ws.Range("M25").Value = 72849428
$ws.Range("M26").Value = 141860762
$ws.Range("M27").Value = 722363250
$ws.Range("M29").Value = 24143510
$ws.Range("M30").Value = "-"
$ws.Range("M31").Value = 1853004
$ws.Range("M32").Value = 35220082
$ws.Range("M33").Value = 1118083
$ws.Range("M34").Value = 52004461
$ws.Range("M35").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("M37").Value = 114339140
$ws.Range("M38").Value = 16400000
$ws.Range("M39").Value = "-"
$ws.Range("M40").Value = 0
$ws.Range("M41").Value = 4100423
$ws.Range("M42").Value = 20500423
$ws.Range("M43").Value = 134839563
$ws.Range("M45").Value = 275000000
$ws.Range("M46").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("M48").Value = -5627478
$ws.Range("M49").Value = 362090
$ws.Range("M50").Value = 27500000
$ws.Range("M51").Value = 0
$ws.Range("M52").Value = "-"
$ws.Range("M53").Value = 0
$ws.Range("M54").Value = "-"
$ws.Range("M55").Value = 0
$ws.Range("M56").Value = 290289075
$ws.Range("M57").Value = 587523687
$ws.Range("M58").Value = 722363250

# Restore column M's width to match the other "publish date" columns
# (D:M originally alternated 29/31-wide columns; M needs the 31-wide style).
$ws.Columns("M").ColumnWidth = 30.1666666666667
